$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for account 004452912 / BRUNO / 47366.32.
# It is row 5 (row 1 is the header: Conta / Nome / Saldo); deleting it
# shifts every following row up by one, matching the diff.
$ws.Rows.Item(5).Delete()
